$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoan_Input")

# -------------------------------------------------------------
# shortname: was the shared text "kar4" -> now numeric 392
# -------------------------------------------------------------
$ws.Range("B3").Value = 392

# -------------------------------------------------------------
# nominalinterestratedefault: 12 -> 1
# -------------------------------------------------------------
$ws.Range("B11").Value = 1

# -------------------------------------------------------------
# New rows 31-42: Ledger/account mapping fields appended below
# the existing data. Copy formatting from rows 10/11 (A uses
# the "label" style, B uses the "value" style) then set values.
# -------------------------------------------------------------
$ws.Range("A10").Copy()
$ws.Range("A31:A42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Copy()
$ws.Range("B31:B42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shared-string table order matches the order cell values are first
# assigned, so write all of column B (the new unique strings) before
# column A (the field-name keys) to reproduce the authored ordering.
$ws.Range("B31").Value = "Cash"
$ws.Range("B32").Value = "Loan portfolio "
$ws.Range("B33").Value = "Interest Receivable "
$ws.Range("B34").Value = "Penalties Receivable "
$ws.Range("B35").Value = "Transfer in Suspence "
$ws.Range("B36").Value = "Fees Receivable"
$ws.Range("B37").Value = "Income from interest"
$ws.Range("B38").Value = "Income from penalties"
$ws.Range("B39").Value = "Income from fees"
$ws.Range("B40").Value = "Income from recovery repayments"
$ws.Range("B41").Value = "Losses Writtenoff "
$ws.Range("B42").Value = "Overpayment Liability"

$ws.Range("A31").Value = "fundsource"
$ws.Range("A32").Value = "loanprotfolio"
$ws.Range("A33").Value = "interestreceivable"
$ws.Range("A34").Value = "penaltiesreceivable"
$ws.Range("A35").Value = "transferinsuspense"
$ws.Range("A36").Value = "feesreceivable"
$ws.Range("A37").Value = "incomefrominterest"
$ws.Range("A38").Value = "incomefrompenalties"
$ws.Range("A39").Value = "incomefromfees"
$ws.Range("A40").Value = "incomefromrecoveryrepayments"
$ws.Range("A41").Value = "loseswrittenoff"
$ws.Range("A42").Value = "overpaymentliability"

# -------------------------------------------------------------
# Update the window/sheet view: active cell moved to B6 and the
# view scrolled down so row 22 is at the top.
# -------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B6").Select()
